$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).AddDays(42604.890833333331)
$ws.Range("B6").Value = "Random"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 65
$ws.Range("I6").Value = 35
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 27
$ws.Range("M6").Value = 73
